$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.744.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +7.01%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.806.06'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.30%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.45'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.44%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9990'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4947'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2774'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +6.35%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06362'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.30%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.799.89'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.69%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.95'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.64%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07085'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.98%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6447'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.51%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.697'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.07%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '81.85'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.44%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.706.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +7.81%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9989'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.03%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007327'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.01%  '

# Row 19
$ws.Range("B19").Value = 'BinanceUSD'
$ws.Range("C19").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9988'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.10%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.21'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.41%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.029.85'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.42%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.609'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.83%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.882'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.62%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.312'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.78%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.74'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.14%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.98'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.885'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +6.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '110.86'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.89%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.387'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.00%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.160'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.79%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08353'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.847'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.54%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04971'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +9.88%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.089'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.84%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.678'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.80%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6681'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9472'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.42%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.625'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +8.27%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.131'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.71%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01596'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.89%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.982'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.60%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.00%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.36'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.07%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4096'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.95%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.219'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.19%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1222'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.01%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05477'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.57%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.122'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.70%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '31.19'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.12%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.304'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.38%  '

# Row 51
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3603'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.46%  '

